# The "Förändrad" (Changed) date in column C was bumped from 2023-09-17
# (serial 45186) to 2023-09-19 (serial 45188) for every data row
# (rows 2 through 453) on the single worksheet in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C453").Value2 = 45188
